$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new rows of day-data (2025-12-27 -> serial 46018)
# Row 54: 四方坪站 (shared string index 4)
$ws.Cells.Item(54, 1).Value = 46018
$ws.Cells.Item(54, 2).Value = "四方坪站"
$ws.Cells.Item(54, 3).Value = 9217.33
$ws.Cells.Item(54, 4).Value = 7725.71
$ws.Cells.Item(54, 5).Value = 3062.87
$ws.Cells.Item(54, 6).Value = 377

# Row 55: 高岭站 (shared string index 5)
$ws.Cells.Item(55, 1).Value = 46018
$ws.Cells.Item(55, 2).Value = "高岭站"
$ws.Cells.Item(55, 3).Value = 4885.12
$ws.Cells.Item(55, 4).Value = 4017.93
$ws.Cells.Item(55, 5).Value = 1295.71
$ws.Cells.Item(55, 6).Value = 182

# Update the selected cell to match the new working position.
$ws.Range("H57").Select() | Out-Null
